# Commit: "some manipulations made with the data as a dictionary"
#
# 1) The X/Y pair that used to read 13 / "15.7" (row 4, index 2 in the
#    dict) is corrected to 13 / "200".
# 2) A new column "z" (a simple range index 1..8) is appended as column D.
# 3) Three more (X, Y) observations are appended as rows 7-9:
#        5, 55, 88
#        6, 70, 300
#        7, 100, 400
# 4) The new "z" column (D2:D9) ends up selected, matching the state the
#    author left the sheet in after adding it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix existing row 4's Y value (stored as text, like its neighbours) ---
$ws.Range("C4").Value = "200"

# --- append the three new observations (rows 7-9) ---
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 55
$ws.Range("C7").Value = 88

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 70
$ws.Range("C8").Value = 300

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 100
$ws.Range("C9").Value = 400

# --- add the "z" column (header + 1..8 index values) ---
$ws.Range("D1").Value = "z"
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 6
$ws.Range("D8").Value = 7
$ws.Range("D9").Value = 8

# --- leave the new column selected, as in the saved workbook ---
$ws.Range("D2:D9").Select() | Out-Null
